# feat: add 2022-Q3 data
#
# Before:
#   "总计"    - overview table; row 2 holds the latest quarter, "2022-Q1"
#   "2022-Q1" - per-fund holding detail for that quarter
#
# After:
#   "总计"    - row 2 now holds the new "2022-Q3" totals, the old
#               "2022-Q1" totals are pushed down to row 3
#   "2022-Q3" - (new) per-fund holding detail for 2022-Q3, inserted
#               right after "总计"
#   "2022-Q1" - per-fund holding detail, unchanged, now the 3rd tab

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" overview sheet: push the existing 2022-Q1 row down to row 3
#    (copy so the A-column style tags along), then fill in row 2 with
#    the new 2022-Q3 totals.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Cells.Item(2, 1).Copy($total.Cells.Item(3, 1))
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q1"
$total.Cells.Item(3, 3).Value = 2
$total.Cells.Item(3, 4).Value = 0.01

$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.06

# ---------------------------------------------------------------------
# 2. Duplicate the existing "2022-Q1" detail sheet so the untouched
#    copy can keep serving as "2022-Q1" after the original is turned
#    into "2022-Q3" below.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($null, $q1)
$q1Copy = $wb.Worksheets.Item($q1.Index + 1)

$q1.Name = "2022-Q3"
$q1Copy.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 3. Turn the original sheet (now named "2022-Q3") into the new
#    quarter's fund holdings. Pull header styling from "总计", which
#    already uses the matching bold/boxed style for its own headers.
# ---------------------------------------------------------------------
$q3 = $q1

$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 2))
$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 3))
$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 4))
$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 5))
$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 6))
$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 7))
$total.Cells.Item(1, 2).Copy($q3.Cells.Item(1, 8))
$total.Cells.Item(2, 1).Copy($q3.Cells.Item(2, 1))
$total.Cells.Item(2, 1).Copy($q3.Cells.Item(3, 1))

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

$q3.Range("B2:G3").NumberFormat = "@"

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "005052"
$q3.Cells.Item(2, 3).Value = "上投摩根标普港股通低波红利指数C"
$q3.Cells.Item(2, 4).Value = "1.37"
$q3.Cells.Item(2, 5).Value = "92.94"
$q3.Cells.Item(2, 6).Value = "2.34"
$q3.Cells.Item(2, 7).Value = "0.0321"
$q3.Cells.Item(2, 8).Value = 7

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "005051"
$q3.Cells.Item(3, 3).Value = "上投摩根标普港股通低波红利指数A"
$q3.Cells.Item(3, 4).Value = "1.36"
$q3.Cells.Item(3, 5).Value = "92.94"
$q3.Cells.Item(3, 6).Value = "2.34"
$q3.Cells.Item(3, 7).Value = "0.0318"
$q3.Cells.Item(3, 8).Value = 7
